$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the existing header cell (G1) into the new header cell (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
